# Add two new columns, I ("I0") and J ("IF"), to the existing table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (row 1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data cells (rows 2-3)
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 5
$ws.Range("I3").Value = 8
$ws.Range("J3").Value = 9

# Match the header formatting used by the rest of row 1 (bold, bordered,
# centered) by copying the style from the existing "IP" header cell (H1).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
